# Auto-generated edit script
# Applies updated market-price-derived values (columns H-N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 349.66666
$ws.Range("I12").Value = 366.6
$ws.Range("K12").Value = 366.6
$ws.Range("M12").Value = -196.6
$ws.Range("H19").Value = 2261.2632
$ws.Range("I19").Value = 2084.9
$ws.Range("J19").Value = 2457.2222
$ws.Range("K19").Value = 2084.9
$ws.Range("L19").Value = 2457.2222
$ws.Range("M19").Value = -1909.9
$ws.Range("N19").Value = -2807.2222
$ws.Range("H33").Value = 5955841
$ws.Range("I33").Value = 6551355
$ws.Range("K33").Value = 6551355
$ws.Range("M33").Value = -6551126
$ws.Range("H74").Value = 4508.143
$ws.Range("I74").Value = 3926.3333
$ws.Range("K74").Value = 3926.3333
$ws.Range("M74").Value = -2990.3333
$ws.Range("H77").Value = 4508.143
$ws.Range("I77").Value = 3926.3333
$ws.Range("K77").Value = 19631.6665
$ws.Range("M77").Value = -14951.6665
$ws.Range("H88").Value = 2363.625
$ws.Range("I88").Value = 1883
$ws.Range("J88").Value = 2844.25
$ws.Range("K88").Value = 1883
$ws.Range("L88").Value = 2844.25
$ws.Range("M88").Value = -1477
$ws.Range("N88").Value = -3656.25
$ws.Range("H91").Value = 2363.625
$ws.Range("I91").Value = 1883
$ws.Range("J91").Value = 2844.25
$ws.Range("K91").Value = 1883
$ws.Range("L91").Value = 2844.25
$ws.Range("M91").Value = -479
$ws.Range("N91").Value = -5652.25
$ws.Range("H99").Value = 1149
$ws.Range("J99").Value = 4500
$ws.Range("L99").Value = 13500
$ws.Range("N99").Value = -16496
$ws.Range("H100").Value = 4255.4
$ws.Range("I100").Value = 3569.5
$ws.Range("K100").Value = 3569.5
$ws.Range("M100").Value = -3028.5
$ws.Range("H101").Value = 903.3333
$ws.Range("I101").Value = 710
$ws.Range("J101").Value = 1290
$ws.Range("K101").Value = 2130
$ws.Range("L101").Value = 3870
$ws.Range("M101").Value = -508
$ws.Range("N101").Value = -7114
$ws.Range("H106").Value = 9998.6
$ws.Range("I106").Value = 9997
$ws.Range("K106").Value = 9997
$ws.Range("M106").Value = -9366
$ws.Range("H132").Value = 62651.35
$ws.Range("I132").Value = 70338.336
$ws.Range("K132").Value = 211015.008
$ws.Range("M132").Value = -208485.008
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 2455.25
$ws.Range("I138").Value = 1140.875
$ws.Range("K138").Value = 3422.625
$ws.Range("M138").Value = 1717.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2312.3
$ws.Range("I61").Value = 2314
$ws.Range("J61").Value = 2297
$ws.Range("K61").Value = 2314
$ws.Range("L61").Value = 2297
$ws.Range("M61").Value = -2102
$ws.Range("N61").Value = -2721
$ws.Range("H74").Value = 15432731
$ws.Range("I74").Value = 9260209
$ws.Range("K74").Value = 9260209
$ws.Range("M74").Value = -9259335
$ws.Range("H77").Value = 15432731
$ws.Range("I77").Value = 9260209
$ws.Range("K77").Value = 46301045
$ws.Range("M77").Value = -46296677
$ws.Range("H136").Value = 2312.3
$ws.Range("I136").Value = 2314
$ws.Range("J136").Value = 2297
$ws.Range("K136").Value = 6942
$ws.Range("L136").Value = 6891
$ws.Range("M136").Value = -4392
$ws.Range("N136").Value = -11991

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 188.88889
$ws.Range("I80").Value = 78
$ws.Range("K80").Value = 78
$ws.Range("M80").Value = 920
$ws.Range("H83").Value = 188.88889
$ws.Range("I83").Value = 78
$ws.Range("K83").Value = 390
$ws.Range("M83").Value = 4602

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3289.842
$ws.Range("I31").Value = 1912
$ws.Range("K31").Value = 1912
$ws.Range("M31").Value = -1617
$ws.Range("H34").Value = 3289.842
$ws.Range("I34").Value = 1912
$ws.Range("K34").Value = 1912
$ws.Range("M34").Value = -1710
$ws.Range("H107").Value = 3106.44
$ws.Range("I107").Value = 2122.2222
$ws.Range("K107").Value = 2122.2222
$ws.Range("M107").Value = -202.2222000000002
$ws.Range("H134").Value = 2732.4375
$ws.Range("I134").Value = 2491
$ws.Range("K134").Value = 7473
$ws.Range("M134").Value = -4938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1104.4
$ws.Range("I7").Value = 934.8
$ws.Range("J7").Value = 1274
$ws.Range("K7").Value = 2804.4
$ws.Range("L7").Value = 3822
$ws.Range("M7").Value = -2692.4
$ws.Range("N7").Value = -4046
$ws.Range("H23").Value = 2048.3914
$ws.Range("I23").Value = 2869.5
$ws.Range("J23").Value = 1416.7693
$ws.Range("K23").Value = 8608.5
$ws.Range("L23").Value = 4250.3079
$ws.Range("M23").Value = -8373.5
$ws.Range("N23").Value = -4720.3079
$ws.Range("H80").Value = 4437.8
$ws.Range("I80").Value = 4150
$ws.Range("K80").Value = 12450
$ws.Range("M80").Value = -11514
$ws.Range("H83").Value = 4437.8
$ws.Range("I83").Value = 4150
$ws.Range("K83").Value = 37350
$ws.Range("M83").Value = -32670
$ws.Range("H104").Value = 4000
$ws.Range("J104").Value = 4000
$ws.Range("L104").Value = 12000
$ws.Range("N104").Value = -17242
$ws.Range("H139").Value = 2910.7693
$ws.Range("I139").Value = 2253.3333
$ws.Range("K139").Value = 6759.999899999999
$ws.Range("M139").Value = -1619.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 36751.332
$ws.Range("J62").Value = 36751.332
$ws.Range("L62").Value = 36751.332
$ws.Range("N62").Value = -38123.332
$ws.Range("H63").Value = 46999.285
$ws.Range("J63").Value = 46999.285
$ws.Range("L63").Value = 46999.285
$ws.Range("N63").Value = -48371.285
$ws.Range("H65").Value = 36751.332
$ws.Range("J65").Value = 36751.332
$ws.Range("L65").Value = 110253.996
$ws.Range("N65").Value = -117117.996
$ws.Range("H66").Value = 46999.285
$ws.Range("J66").Value = 46999.285
$ws.Range("L66").Value = 140997.855
$ws.Range("N66").Value = -147861.855
$ws.Range("H69").Value = 58749.75
$ws.Range("J69").Value = 58749.75
$ws.Range("L69").Value = 58749.75
$ws.Range("N69").Value = -60247.75
$ws.Range("H72").Value = 58749.75
$ws.Range("J72").Value = 58749.75
$ws.Range("L72").Value = 176249.25
$ws.Range("N72").Value = -183737.25
$ws.Range("H122").Value = 3235.7917
$ws.Range("I122").Value = 3364.7368
$ws.Range("K122").Value = 10094.2104
$ws.Range("M122").Value = -7644.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 567.4857
$ws.Range("I55").Value = 639.73334
$ws.Range("J55").Value = 513.3
$ws.Range("K55").Value = 639.73334
$ws.Range("L55").Value = 513.3
$ws.Range("M55").Value = -466.73334
$ws.Range("N55").Value = -859.3
$ws.Range("H132").Value = 2665.5557
$ws.Range("I132").Value = 2427.2856
$ws.Range("K132").Value = 7281.8568
$ws.Range("M132").Value = -4751.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 19000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 19000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -19224
$ws.Range("H26").Value = 8548.25
$ws.Range("I26").Value = 3064.6667
$ws.Range("J26").Value = 24999
$ws.Range("K26").Value = 3064.6667
$ws.Range("L26").Value = 24999
$ws.Range("M26").Value = -2771.6667
$ws.Range("N26").Value = -25585
$ws.Range("H32").Value = 11342.833
$ws.Range("I32").Value = 11000
$ws.Range("J32").Value = 11514.25
$ws.Range("K32").Value = 11000
$ws.Range("L32").Value = 11514.25
$ws.Range("M32").Value = -10683
$ws.Range("N32").Value = -12148.25

Write-Host "Applied 202 cell updates and 2 cell clears."
